$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.312.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.691.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.68%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'217.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.19%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5398"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.77%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2734"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.24%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06454"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.41%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07675"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.89%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'1.719.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.38%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'4.543"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.5792"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.19%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.000008397"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.66%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +3.46%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.374.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.17%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.910"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.25%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.12%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.05%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'190.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.261"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.91%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'149.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.41%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.99%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.869"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'15.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.30%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.06292"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.374"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.26%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.326"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.604"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.30%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.582"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.59%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.678"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.00%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.07%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.6178"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.47%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.770"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.91%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.02%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.111.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.02%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.21%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8829"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'101.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.53%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.843.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.72%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'57.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("B47").Value = "'Frax"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.006"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.18%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'8.150"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.73%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.05283"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.37%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.4299"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.18%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.07%  "
$ws.Range("E51").Style = "Normal"
